# Automatische test-sync: 2025-06-19 19:05:30
# Append a new "Afmelding nieuwsbrief" log entry to the Logs sheet,
# extend the conditional formatting ranges to cover the new row,
# and update the Dashboard's "Afmelding" count accordingly.

$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

$newRow = 36

$logs.Cells.Item($newRow, 1).Value = "Afmelding nieuwsbrief"
$logs.Cells.Item($newRow, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newRow, 3).Value = "Graag afmelden voor de nieuwsbrief. Dank u."
$logs.Cells.Item($newRow, 4).Value = "Afmelding"
$logs.Cells.Item($newRow, 6).Value = "2025-06-19 19:05:25"
$logs.Cells.Item($newRow, 7).Value = "Nee"

# Extend the conditional formatting sqref ranges (D2:D35 -> D2:D36, G2:G35 -> G2:G36)
$dFormatConditions = $logs.Range("D2:D35").FormatConditions
$dFormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D36"))

$gFormatConditions = $logs.Range("G2:G35").FormatConditions
$gFormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G36"))

# Update the Dashboard summary count for "Afmelding"
$dashboard.Range("B2").Value = 9
